$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted right above the current
# row 74. Insert a whole row there so every subsequent record (old rows
# 74-111) shifts down by one and keeps its formatting (row 112 is created
# at the bottom of the used range as a result).
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new record. The static
# descriptive columns (market/region/product taxonomy) repeat the same
# values used throughout this block of rows.
$ws.Cells.Item(74, 1).Value  = 5
$ws.Cells.Item(74, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(74, 3).Value  = "Maule"
$ws.Cells.Item(74, 4).Value  = 44609
$ws.Cells.Item(74, 5).Value  = 7
$ws.Cells.Item(74, 6).Value  = "Fruta"
$ws.Cells.Item(74, 7).Value  = 100108
$ws.Cells.Item(74, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(74, 9).Value  = 100108002
$ws.Cells.Item(74, 10).Value = "Mango"
$ws.Cells.Item(74, 11).Value = "Sin especificar"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 150
$ws.Cells.Item(74, 14).Value = 7000
$ws.Cells.Item(74, 15).Value = 7000
$ws.Cells.Item(74, 16).Value = 7000
$ws.Cells.Item(74, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(74, 18).Value = "Perú"
$ws.Cells.Item(74, 19).Value = 1750
$ws.Cells.Item(74, 20).Value = 4
